$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fzd8"
$ws.Range("C2").Value = "Ckap4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.027778333333333
$ws.Range("H2").Value = 3.083335
$ws.Range("I2").Value = 0.08020467841353289
$ws.Range("J2").Value = 0.08020467841353289
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.44960333333333
$ws.Range("N2").Value = 43.34881
$ws.Range("O2").Value = 0.4336357442982187
$ws.Range("P2").Value = 0.4336357442982187
$ws.Range("Q2").Value = 14.85098923126111
$ws.Range("R2").Value = 133.65890308135
$ws.Range("S2").Value = 0.0347796154200516
$ws.Range("T2").Value = 0.0347796154200516
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fzd8"
$ws.Range("C3").Value = "Ckap4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.027778333333333
$ws.Range("H3").Value = 3.083335
$ws.Range("I3").Value = 0.08020467841353289
$ws.Range("J3").Value = 0.08020467841353289
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.96339
$ws.Range("N3").Value = 44.89017
$ws.Range("O3").Value = 0.4490545941082019
$ws.Range("P3").Value = 0.4490545941082019
$ws.Range("Q3").Value = 15.37904803521666
$ws.Range("R3").Value = 138.41143231695
$ws.Range("S3").Value = 0.03601627931056788
$ws.Range("T3").Value = 0.03601627931056788
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fzd8"
$ws.Range("C4").Value = "Ckap4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.027778333333333
$ws.Range("H4").Value = 3.083335
$ws.Range("I4").Value = 0.08020467841353289
$ws.Range("J4").Value = 0.08020467841353289
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.908990666666667
$ws.Range("N4").Value = 11.726972
$ws.Range("O4").Value = 0.1173096615935794
$ws.Range("P4").Value = 0.1173096615935794
$ws.Range("Q4").Value = 4.017575912402222
$ws.Range("R4").Value = 36.15818321162
$ws.Range("S4").Value = 0.009408783682913403
$ws.Range("T4").Value = 0.009408783682913403
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fzd8"
$ws.Range("C5").Value = "Ckap4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.273511666666667
$ws.Range("H5").Value = 21.820535
$ws.Range("I5").Value = 0.567602609669802
$ws.Range("J5").Value = 0.567602609669802
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 14.44960333333333
$ws.Range("N5").Value = 43.34881
$ws.Range("O5").Value = 0.4336357442982187
$ws.Range("P5").Value = 0.4336357442982187
$ws.Range("Q5").Value = 105.0993584237056
$ws.Range("R5").Value = 945.89422581335
$ws.Range("S5").Value = 0.2461327801097759
$ws.Range("T5").Value = 0.2461327801097759
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fzd8"
$ws.Range("C6").Value = "Ckap4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.273511666666667
$ws.Range("H6").Value = 21.820535
$ws.Range("I6").Value = 0.567602609669802
$ws.Range("J6").Value = 0.567602609669802
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.96339
$ws.Range("N6").Value = 44.89017
$ws.Range("O6").Value = 0.4490545941082019
$ws.Range("P6").Value = 0.4490545941082019
$ws.Range("Q6").Value = 108.8363917378833
$ws.Range("R6").Value = 979.5275256409499
$ws.Range("S6").Value = 0.2548845595000291
$ws.Range("T6").Value = 0.2548845595000291
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fzd8"
$ws.Range("C7").Value = "Ckap4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.273511666666667
$ws.Range("H7").Value = 21.820535
$ws.Range("I7").Value = 0.567602609669802
$ws.Range("J7").Value = 0.567602609669802
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.908990666666667
$ws.Range("N7").Value = 11.726972
$ws.Range("O7").Value = 0.1173096615935794
$ws.Range("P7").Value = 0.1173096615935794
$ws.Range("Q7").Value = 28.43208921889111
$ws.Range("R7").Value = 255.88880297002
$ws.Range("S7").Value = 0.066585270059997
$ws.Range("T7").Value = 0.066585270059997
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fzd8"
$ws.Range("C8").Value = "Ckap4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.513153666666667
$ws.Range("H8").Value = 13.539461
$ws.Range("I8").Value = 0.3521927119166651
$ws.Range("J8").Value = 0.3521927119166651
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.44960333333333
$ws.Range("N8").Value = 43.34881
$ws.Range("O8").Value = 0.4336357442982187
$ws.Range("P8").Value = 0.4336357442982187
$ws.Range("Q8").Value = 65.21328026571223
$ws.Range("R8").Value = 586.9195223914099
$ws.Range("S8").Value = 0.1527233487683912
$ws.Range("T8").Value = 0.1527233487683912
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fzd8"
$ws.Range("C9").Value = "Ckap4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.513153666666667
$ws.Range("H9").Value = 13.539461
$ws.Range("I9").Value = 0.3521927119166651
$ws.Range("J9").Value = 0.3521927119166651
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 14.96339
$ws.Range("N9").Value = 44.89017
$ws.Range("O9").Value = 0.4490545941082019
$ws.Range("P9").Value = 0.4490545941082019
$ws.Range("Q9").Value = 67.53207844426333
$ws.Range("R9").Value = 607.7887059983699
$ws.Range("S9").Value = 0.1581537552976049
$ws.Range("T9").Value = 0.1581537552976049
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fzd8"
$ws.Range("C10").Value = "Ckap4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.513153666666667
$ws.Range("H10").Value = 13.539461
$ws.Range("I10").Value = 0.3521927119166651
$ws.Range("J10").Value = 0.3521927119166651
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.908990666666667
$ws.Range("N10").Value = 11.726972
$ws.Range("O10").Value = 0.1173096615935794
$ws.Range("P10").Value = 0.1173096615935794
$ws.Range("Q10").Value = 17.64187556023244
$ws.Range("R10").Value = 158.776880042092
$ws.Range("S10").Value = 0.04131560785066898
$ws.Range("T10").Value = 0.04131560785066897